$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking": Right marks 4 -> 5, Wrong marks -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 "Total": Right total 72 -> 90, Wrong total -7 -> -8.4, Max string updated
$ws.Range("B12").Value = 90
$ws.Range("C12").Value = -8.4
$ws.Range("E12").Value = "81.6/140"
